$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 128; this shifts the existing rows 128-140
# down to 129-141 (Excel also duplicates the row-128 cell styles, e.g. the
# date style on column D, onto the freshly inserted row).
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with the new weekly price record.
$ws.Cells.Item(128, 1).Value = 10
$ws.Cells.Item(128, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(128, 3).Value = "La Araucanía"
$ws.Cells.Item(128, 4).Value = 44746
$ws.Cells.Item(128, 5).Value = 9
$ws.Cells.Item(128, 6).Value = 100112031
$ws.Cells.Item(128, 7).Value = "Poroto verde"
$ws.Cells.Item(128, 8).Value = "Sin especificar"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 125
$ws.Cells.Item(128, 11).Value = 35000
$ws.Cells.Item(128, 12).Value = 35000
$ws.Cells.Item(128, 13).Value = 35000
$ws.Cells.Item(128, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(128, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(128, 16).Value = 1400
$ws.Cells.Item(128, 17).Value = 25
$ws.Cells.Item(128, 18).Value = "Hortaliza"
